# Fill in the computed Retention ratio and ARA (Answer Recall) values
# that were left blank in the summary table.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-ComputedValue($rowIndex, $value) {
    $cell = $t.Rows.Item($rowIndex).Cells.Item(2)
    $r = $cell.Range
    $r.Text = $value
    $r.Font.Bold = $true
    $r.Font.Size = 12
    $r.Font.SizeBi = 12
}

# RETENTION section - "Ratio" row
Set-ComputedValue 24 "0.8333"

# QUESTION AND ANSWER TASK section - Answer Recall rows
Set-ComputedValue 44 "0.375"
Set-ComputedValue 45 "0.125"
Set-ComputedValue 46 "0.25"
